$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B: replace per-student spending figures with new (smaller) integer values ---
$ws.Range("B2").Value  = 3770
$ws.Range("B3").Value  = 8482
$ws.Range("B4").Value  = 3703
$ws.Range("B5").Value  = 7426
$ws.Range("B6").Value  = 8153
$ws.Range("B7").Value  = 2617
$ws.Range("B8").Value  = 7749
$ws.Range("B9").Value  = 5616
$ws.Range("B10").Value = 6921
$ws.Range("B11").Value = 4232
$ws.Range("B12").Value = 6571
$ws.Range("B13").Value = 7496
$ws.Range("B14").Value = 6666
$ws.Range("B15").Value = 10603
$ws.Range("B16").Value = 9572
$ws.Range("B17").Value = 10807
$ws.Range("B18").Value = 5982
$ws.Range("B19").Value = 7584
$ws.Range("B20").Value = 6092
$ws.Range("B21").Value = 3618
$ws.Range("B22").Value = 9251
$ws.Range("B23").Value = 4518
$ws.Range("B24").Value = 2963
$ws.Range("B25").Value = 13528

# --- Column E: RANK formulas now look at $B$2:$B$25 instead of $B$2:$B$26 ---
# E2 is a standalone (non-shared) formula.
$ws.Range("E2").Formula = '=RANK(B2,$B$2:$B$25,0)'

# E3:E25 used to be one shared-formula block (E3:E26). Row 17 got a one-off
# override (`=E19`) typed over it, so we rebuild the block around that gap.
$ws.Range("E3:E16").Formula  = '=RANK(B3,$B$2:$B$25,0)'
$ws.Range("E18:E25").Formula = '=RANK(B18,$B$2:$B$25,0)'
$ws.Range("E17").Formula = '=E19'

# --- Row 26 ("Maryland" summary row): spending + both rank columns cleared ---
$ws.Range("B26").ClearContents()
$ws.Range("D26").ClearContents()
$ws.Range("E26").ClearContents()

# --- Selection moves to E15 ---
$ws.Range("E15").Select()
